# Apply "Natmi following Dr Hou advice" update to Omg-Lingo1 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing values
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Omg"
$ws.Cells.Item(2, 3).Value = "Lingo1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 4.523724666666666
$ws.Cells.Item(2, 8).Value = 13.571174
$ws.Cells.Item(2, 9).Value = 0.4806607624766543
$ws.Cells.Item(2, 10).Value = 0.4806607624766543
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.9237416666666666
$ws.Cells.Item(2, 14).Value = 2.771225
$ws.Cells.Item(2, 15).Value = 0.9776944302049534
$ws.Cells.Item(2, 16).Value = 0.9776944302049534
$ws.Cells.Item(2, 17).Value = 4.178752963127777
$ws.Cells.Item(2, 18).Value = 37.60877666814999
$ws.Cells.Item(2, 19).Value = 0.469939350291491
$ws.Cells.Item(2, 20).Value = 0.469939350291491

# Row 3: update existing values
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Omg"
$ws.Cells.Item(3, 3).Value = "Lingo1"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 4.523724666666666
$ws.Cells.Item(3, 8).Value = 13.571174
$ws.Cells.Item(3, 9).Value = 0.4806607624766543
$ws.Cells.Item(3, 10).Value = 0.4806607624766543
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.02107466666666667
$ws.Cells.Item(3, 14).Value = 0.063224
$ws.Cells.Item(3, 15).Value = 0.02230556979504659
$ws.Cells.Item(3, 16).Value = 0.02230556979504659
$ws.Cells.Item(3, 17).Value = 0.09533598944177778
$ws.Cells.Item(3, 18).Value = 0.858023904976
$ws.Cells.Item(3, 19).Value = 0.01072141218516332
$ws.Cells.Item(3, 20).Value = 0.01072141218516332

# Row 4: update existing values
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Omg"
$ws.Cells.Item(4, 3).Value = "Lingo1"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.109481
$ws.Cells.Item(4, 8).Value = 6.328443
$ws.Cells.Item(4, 9).Value = 0.2241393587371326
$ws.Cells.Item(4, 10).Value = 0.2241393587371326
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.9237416666666666
$ws.Cells.Item(4, 14).Value = 2.771225
$ws.Cells.Item(4, 15).Value = 0.9776944302049534
$ws.Cells.Item(4, 16).Value = 0.9776944302049534
$ws.Cells.Item(4, 17).Value = 1.948615494741667
$ws.Cells.Item(4, 18).Value = 17.537539452675
$ws.Cells.Item(4, 19).Value = 0.2191398026270045
$ws.Cells.Item(4, 20).Value = 0.2191398026270045

# Row 5: new row
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Omg"
$ws.Cells.Item(5, 3).Value = "Lingo1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.109481
$ws.Cells.Item(5, 8).Value = 6.328443
$ws.Cells.Item(5, 9).Value = 0.2241393587371326
$ws.Cells.Item(5, 10).Value = 0.2241393587371326
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.02107466666666667
$ws.Cells.Item(5, 14).Value = 0.063224
$ws.Cells.Item(5, 15).Value = 0.02230556979504659
$ws.Cells.Item(5, 16).Value = 0.02230556979504659
$ws.Cells.Item(5, 17).Value = 0.04445660891466668
$ws.Cells.Item(5, 18).Value = 0.400109480232
$ws.Cells.Item(5, 19).Value = 0.004999556110128095
$ws.Cells.Item(5, 20).Value = 0.004999556110128095

# Row 6: new row
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Omg"
$ws.Cells.Item(6, 3).Value = "Lingo1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.778265
$ws.Cells.Item(6, 8).Value = 8.334795
$ws.Cells.Item(6, 9).Value = 0.2951998787862131
$ws.Cells.Item(6, 10).Value = 0.2951998787862131
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.9237416666666666
$ws.Cells.Item(6, 14).Value = 2.771225
$ws.Cells.Item(6, 15).Value = 0.9776944302049534
$ws.Cells.Item(6, 16).Value = 0.9776944302049534
$ws.Cells.Item(6, 17).Value = 2.566399141541666
$ws.Cells.Item(6, 18).Value = 23.097592273875
$ws.Cells.Item(6, 19).Value = 0.2886152772864579
$ws.Cells.Item(6, 20).Value = 0.2886152772864579

# Row 7: new row
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Omg"
$ws.Cells.Item(7, 3).Value = "Lingo1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.778265
$ws.Cells.Item(7, 8).Value = 8.334795
$ws.Cells.Item(7, 9).Value = 0.2951998787862131
$ws.Cells.Item(7, 10).Value = 0.2951998787862131
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.02107466666666667
$ws.Cells.Item(7, 14).Value = 0.063224
$ws.Cells.Item(7, 15).Value = 0.02230556979504659
$ws.Cells.Item(7, 16).Value = 0.02230556979504659
$ws.Cells.Item(7, 17).Value = 0.05855100878666666
$ws.Cells.Item(7, 18).Value = 0.5269590790799999
$ws.Cells.Item(7, 19).Value = 0.006584601499755168
$ws.Cells.Item(7, 20).Value = 0.006584601499755168
